$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: remove the out-of-range forecast values (C2, E2) produced by the
# naive component forecaster bug.
$ws.Range("C2").ClearContents()
$ws.Range("E2").ClearContents()

# Row 3: remove the out-of-range forecast value (C3) and correct E3.
$ws.Range("C3").ClearContents()
$ws.Range("E3").Value = 6.772115316529947

# Row 4: corrected value after the forecaster fix.
$ws.Range("C4").Value = -7.266312015249799

# Row 6: corrected value after the forecaster fix.
$ws.Range("C6").Value = 9.469137444079955

# Row 7: corrected values after the forecaster fix.
$ws.Range("C7").Value = 3.358206407534969
$ws.Range("E7").Value = 5.745831525574463

# Row 9: corrected values after the forecaster fix.
$ws.Range("C9").Value = 3.901355411819685
$ws.Range("E9").Value = 3.690459963535031

# Row 12: corrected value after the forecaster fix.
$ws.Range("C12").Value = 5.246209615995689

# Row 13: corrected values after the forecaster fix.
$ws.Range("C13").Value = 4.862559663742938
$ws.Range("E13").Value = 4.636196713604379

# Row 14: corrected value after the forecaster fix.
$ws.Range("C14").Value = 2.76474001115945

# Row 15: corrected value after the forecaster fix.
$ws.Range("C15").Value = -7.260793671746447

# Row 16: corrected value after the forecaster fix.
$ws.Range("C16").Value = 4.097586525396246

# Row 17: corrected value after the forecaster fix.
$ws.Range("C17").Value = 7.824284864703768

# Row 18: corrected value after the forecaster fix.
$ws.Range("C18").Value = -1.245022353133318
